$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (65) into the new row (66)
$ws.Range("A65:F65").Copy()
$ws.Range("A66:F66").PasteSpecial(-4122)

# Populate the new row's values
$ws.Range("A66").Value = "L_ERSTT_10"
$ws.Range("B66").Value = "Q_EUROSTAT"
$ws.Range("C66").Value = "Rate der erheblichen materiellen und sozialen Deprivation"
$ws.Range("D66").Value = "Severe material and social deprivation rate"
$ws.Range("E66").Value = ""
$ws.Range("F66").Value = "https://ec.europa.eu/eurostat/databrowser/view/ILC_MDSD11__custom_3696252/default/table?lang=de"
